$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Fix up the boolean driver_info ("safe_driving" / "is_primary_driver")
#    cells on sheet 2 so they are formulas (TRUE()/FALSE()) instead of
#    bare boolean literals - same values, but now formula-backed.
# ---------------------------------------------------------------------
$driverInfo = $wb.Worksheets.Item("driver_info")

$driverInfo.Range("E2").Formula = "=FALSE()"
$driverInfo.Range("F2").Formula = "=TRUE()"

$driverInfo.Range("E3").Formula = "=FALSE()"
$driverInfo.Range("F3").Formula = "=FALSE()"

$driverInfo.Range("E4").Formula = "=TRUE()"
$driverInfo.Range("F4").Formula = "=TRUE()"

$driverInfo.Range("E5").Formula = "=FALSE()"
$driverInfo.Range("F5").Formula = "=FALSE()"

$driverInfo.Range("E6").Formula = "=TRUE()"
$driverInfo.Range("F6").Formula = "=TRUE()"

$driverInfo.Range("E7").Formula = "=FALSE()"
$driverInfo.Range("F7").Formula = "=TRUE()"

# ---------------------------------------------------------------------
# 2. Add the new "rating_inputs_simple" worksheet at the end of the
#    workbook and populate it with the sample rating-input rows.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ratingInputs = $wb.Worksheets.Add($null, $lastSheet)
$ratingInputs.Name = "rating_inputs_simple"

$ratingInputs.Range("A1").Value = "age"
$ratingInputs.Range("B1").Value = "safe_driving"
$ratingInputs.Range("C1").Value = "credit_tier"

$ratingInputs.Range("A2").Value = 1
$ratingInputs.Range("B2").Value = $true
$ratingInputs.Range("C2").Value = "C1"

$ratingInputs.Range("A3").Value = 18
$ratingInputs.Range("B3").Value = $true
$ratingInputs.Range("C3").Value = "B1"

$ratingInputs.Range("A4").Value = 25
$ratingInputs.Range("B4").Value = "null"
$ratingInputs.Range("C4").Value = "D1"

$ratingInputs.Range("A5").Value = 470
$ratingInputs.Range("B5").Value = $false
$ratingInputs.Range("C5").Value = "E1"

$ratingInputs.Range("A6").Value = -1
$ratingInputs.Range("B6").Value = "Missing"
$ratingInputs.Range("C6").Value = "X1"

# ---------------------------------------------------------------------
# 3. Bump every sheet's zoom level from 120% to 130%, keeping gridlines
#    visible (their display state doesn't otherwise change).
# ---------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Activate()
    $win = $excel.ActiveWindow
    $win.DisplayGridlines = $true
    $win.Zoom = 130
}

# ---------------------------------------------------------------------
# 4. The new sheet becomes the active / selected tab.
# ---------------------------------------------------------------------
$ratingInputs.Activate()
$ratingInputs.Range("B6").Select()
